$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (feature #4): rename the task and add a note explaining the scope was
# reduced to "view only" for now.
$ws.Range("B5").Value = "Xem nội dung tin rao vặt."
$ws.Range("F5").Value = "Chỉ mới xem được tin rao vặt thường"

# The longer note now wraps onto two lines, so the row grows taller.
$ws.Rows.Item(5).RowHeight = 30

# Move the active selection from F16 to D14, matching the editor's cursor
# position when this edit was made.
$ws.Range("D14").Select() | Out-Null
